# Rename the AHB-Diff column headers so the "_old"/"_new" suffixes reflect
# the two format versions being compared (FV2210 vs FV2304), then wrap the
# sheet's data range in a real Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells in row 1 -----------------------------------
# Columns A:J carry the "_old" suffix -> "_FV2210"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2210')
}

# Columns L:U carry the "_new" suffix -> "_FV2304" (column K is just "diff")
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2304')
}

# --- 2. Turn the used range into an Excel Table -----------------------------
$tableRange = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row (top row) ------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
